$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Remove the stray "_GoBack" bookmark from the "General information"
#    heading paragraph (it gets re-created later, at the end of the new
#    "{% else -%}" paragraph).
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 2. After the "Results" heading, add a new body paragraph that opens
#    the Jinja "if" block:  {% if project.vulnComponentsCount -%}
# ---------------------------------------------------------------------
$resultsPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim() -eq "Results") {
        $resultsPara = $p
    }
}
$resultsPara.Range.InsertParagraphAfter()
$ifPara = $resultsPara.Next()
$ifPara.Style = $d.Styles("Normal")
$ifPara.Range.LanguageID = "en-US"
$ifPara.Range.InsertAfter("{% if project.vulnComponentsCount -%}")
$ifPara.Range.LanguageID = "en-US"

# ---------------------------------------------------------------------
# 3. The three table rows that used to be marked "cantSplit" (can't
#    break across a page) are allowed to split again.
# ---------------------------------------------------------------------
$tbl = $d.Tables(1)
$tbl.Rows(2).AllowBreakAcrossPages = $true
$tbl.Rows(3).AllowBreakAcrossPages = $true
$tbl.Rows(4).AllowBreakAcrossPages = $true

# ---------------------------------------------------------------------
# 4. Replace the single trailing empty paragraph (right after the
#    table, before the section break) with three new paragraphs that
#    close out the template:
#       {% else -%}
#       No vulnerabilities were found.
#       {%- endif %}
#    The "_GoBack" bookmark is re-added, collapsed, right at the end
#    of the "{% else -%}" paragraph (after its text, before the
#    paragraph mark) -- exactly where it used to live.
# ---------------------------------------------------------------------
$trailingPara = $d.Paragraphs($d.Paragraphs.Count)

$trailingPara.Range.LanguageID = "en-US"
$trailingPara.Range.InsertAfter("{% else -%}")
$trailingPara.Range.LanguageID = "en-US"

# Re-add the "_GoBack" bookmark at the very end of the "{% else -%}"
# paragraph, collapsed (zero length), right before its paragraph mark.
# A direct collapsed Range at that exact boundary position is mis-
# resolved by the engine, so a one-character placeholder is used to
# hold the spot open while the bookmark is inserted, then removed.
$trailingPara.Range.InsertAfter("X")
$full = $trailingPara.Range
$xPos = $full.End - 2
$markRange = $d.Range($xPos, $xPos)
$d.Bookmarks.Add("_GoBack", $markRange)
$d.Range($xPos, $xPos + 1).Delete()

$trailingPara.Range.InsertParagraphAfter()
$noVulnPara = $trailingPara.Next()
$noVulnPara.Range.LanguageID = "en-US"
$noVulnPara.Range.InsertAfter("No vulnerabilities were found.")
$noVulnPara.Range.LanguageID = "en-US"

$noVulnPara.Range.InsertParagraphAfter()
$endifPara = $noVulnPara.Next()
$endifPara.Range.LanguageID = "en-US"
$endifPara.Range.InsertAfter("{%- endif %}")
$endifPara.Range.LanguageID = "en-US"
